# Update the "Metadata" sheet with the new IG identity/version/date/publisher.
$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/related-issue"
$wsMeta.Range("B3").Value = "8.0.0"
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# Update the "Elements" sheet: the StructureDefinition base/reference URLs move
# from ibm.com to linuxforhealth.org, and the stale constraint text duplicated
# onto the root "Extension" row is cleared (it legitimately belongs only to the
# "Extension.extension" and "Extension.value[x]" rows).
$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("AI2").Value = ""
$wsElem.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/related-issue"
$wsElem.Range("J6").Value = "Reference(http://linuxforhealth.org/fhir/cdm/StructureDefinition/care-gap-detected-issue)`n"
